$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = "INCLUDE"
$ws.Range("D2").Value = "DUE"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "The abstract is not provided. As per the abstract completeness rule, the article is included because a proper evaluation of the study type, device relevance, and clinical data cannot be performed based on the missing information. A full-text review is required for a complete assessment."

# Row 3
$ws.Range("C3").Value = "EXCLUDE"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "Literature contains insufficient information to undertake a scientific analysis about device performance."
$ws.Range("F3").Value = "The abstract provided is empty and contains no information. Therefore, it is not possible to conduct a scientific analysis or determine the relevance of the study, its design, or its outcomes."

# Row 4
$ws.Range("C4").Value = "INCLUDE"
$ws.Range("D4").Value = "DUE"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "The abstract describes a clinical study investigating the correlation between subthalamic nucleus (STN) local field potentials (LFPs) and motor impairment in Parkinson's disease patients. The study utilizes a sensing-enabled implantable pulse generator, which directly aligns with the features and intended use of the subject device (Percept PC Neurostimulator). The abstract presents relevant clinical performance data on the device's sensing capabilities, and no exclusion criteria are met."

# Row 5
$ws.Range("C5").Value = "EXCLUDE"
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "Pre-clinical, animal, cadaver, In-vitro, benchtop, biomechanical, or other non-clinical study"
$ws.Range("F5").Value = "The abstract explicitly states that the study is a biomechanical evaluation conducted on a calf spine model. This is a pre-clinical, animal study. This study type is a direct match for the exclusion criteria."
